# Auto-update draw results: append the 2025-12-04 Pick 4 draw as row 79,
# mirroring the existing rows (2-78) whose A/B/C/D/E cells are all stored
# as plain text using the workbook's one-and-only (default) cell style.
#
# Columns A (date, e.g. "2025-12-04") and C (phase, e.g. "251204") look
# like a date / a number respectively, so a plain
#   $ws.Range("A79").Value = "2025-12-04"
# assignment would make Excel auto-convert them into a real date serial
# / numeric value (and any apostrophe-prefix / NumberFormat="@" trick to
# force text leaves a stray "quote prefix" style on the cell, unlike the
# untouched default style the rest of the sheet uses).
#
# Instead, write them as a formula that evaluates to a text literal, then
# copy/paste-special-values over itself: the result is a genuine stored
# string (not a formula, not a number/date) in the default style - an
# exact match for how the pre-existing rows are encoded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 79
$xlPasteValues = -4163

$date = "2025-12-04"
$game = "Pick 4"
$phase = "251204"
$result = "9-3-2-4"
$insertedAt = "2025-12-04T21:44:11.961+04:00"

function Set-TextCell($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.Formula = '="' + $text + '"'
    $rng.Copy()
    $rng.PasteSpecial($xlPasteValues)
}

# A79 - date-like text ("2025-12-04" would otherwise become a date serial).
Set-TextCell "A$newRow" $date

# B79 - plain text, no ambiguity; a direct value assignment is fine.
$ws.Range("B$newRow").Value = $game

# C79 - numeric-looking text ("251204" would otherwise become a number).
Set-TextCell "C$newRow" $phase

# D79 - plain text (contains dashes, not parsed as a number/date).
$ws.Range("D$newRow").Value = $result

# E79 - ISO-ish timestamp text, not parsed as a date by Excel's inference.
$ws.Range("E$newRow").Value = $insertedAt
